$d = $word.ActiveDocument

# --- Locate the final paragraph of the body (currently the empty
#     "0F672616" paragraph that sits right before the sectPr). ---
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $pLast.Range

# Insert four new empty paragraphs right before it - this pushes the
# bookmark-carrying slot down and gives us room for the new content
# without disturbing the trailing empty paragraph's identity.
$rLast.InsertParagraphBefore()
$rLast.InsertParagraphBefore()
$rLast.InsertParagraphBefore()
$rLast.InsertParagraphBefore()

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Paragraph with "...7e9ac5de74s1d": strip the _GoBack bookmark
#     from it (it moves down to the new trailing paragraph below). ---
$pCode = $d.Paragraphs.Item(22)
$pCode.Range.InsertXML("<w:p $wNs>" +
  "<w:r><w:t>Privil" + [char]0x00E9 + "gi" + [char]0x00E9 + [char]0x00A0 + ": 7</w:t></w:r>" +
  "<w:r><w:t>e</w:t></w:r>" +
  "<w:r><w:t>9ac5de74s1d</w:t></w:r>" +
  "</w:p>")

# --- New "Switchs : ..." paragraph ---
$pSwitchs = $d.Paragraphs.Item(24)
$pSwitchs.Range.InsertXML("<w:p $wNs>" +
  "<w:proofErr w:type='spellStart'/>" +
  "<w:r><w:t>Switchs</w:t></w:r>" +
  "<w:proofErr w:type='spellEnd'/>" +
  "<w:r><w:t xml:space='preserve'> : </w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'>aghjklomp+2lettres du </w:t></w:r>" +
  "<w:proofErr w:type='spellStart'/>" +
  "<w:r><w:t>jeu+numero</w:t></w:r>" +
  "<w:proofErr w:type='spellEnd'/>" +
  "<w:r><w:t xml:space='preserve'> d switch</w:t></w:r>" +
  "</w:p>")

# --- New "Ex : SwitchPU2" paragraph ---
$pEx = $d.Paragraphs.Item(25)
$pEx.Range.InsertXML("<w:p $wNs><w:r><w:t>Ex : SwitchPU2</w:t></w:r></w:p>")

# --- New "MDP : ..." paragraph (starts with a tab) ---
$pMdp = $d.Paragraphs.Item(26)
$pMdp.Range.InsertXML("<w:p $wNs>" +
  "<w:r><w:tab/><w:t xml:space='preserve'>MDP : </w:t></w:r>" +
  "<w:r><w:t>aghjklomp</w:t></w:r>" +
  "<w:r><w:t>PU2</w:t></w:r>" +
  "</w:p>")

# --- Final (previously-last) paragraph now carries the _GoBack bookmark alone ---
$pBookmark = $d.Paragraphs.Item(27)
$pBookmark.Range.InsertXML("<w:p $wNs>" +
  "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
  "<w:bookmarkEnd w:id='0'/>" +
  "</w:p>")
